# Update gh-pages output data (江西-漫展信息.xlsx)
# Applies updated "想去人数" (F) and "最低票价" (G) figures to the
# "展览" (Exhibition) sheet and the combined "全部类型" (All types) sheet.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (CellAddress, NewValue) updates
$updates = @{
    "展览" = @(
        @{ Cell = "G2";  Value = 50 }
        @{ Cell = "F3";  Value = 1881 }
        @{ Cell = "G3";  Value = 65 }
        @{ Cell = "G4";  Value = 68 }
        @{ Cell = "G5";  Value = 25 }
        @{ Cell = "F6";  Value = 850 }
        @{ Cell = "G6";  Value = 60 }
        @{ Cell = "F16"; Value = 4430 }
        @{ Cell = "F19"; Value = 482 }
        @{ Cell = "F22"; Value = 13 }
        @{ Cell = "F23"; Value = 1108 }
        @{ Cell = "F24"; Value = 2002 }
        @{ Cell = "F29"; Value = 2118 }
        @{ Cell = "F30"; Value = 77 }
        @{ Cell = "F33"; Value = 150 }
        @{ Cell = "F35"; Value = 35 }
    )
    "全部类型" = @(
        @{ Cell = "G2";  Value = 50 }
        @{ Cell = "F3";  Value = 1881 }
        @{ Cell = "G3";  Value = 65 }
        @{ Cell = "G4";  Value = 68 }
        @{ Cell = "G5";  Value = 25 }
        @{ Cell = "F6";  Value = 850 }
        @{ Cell = "G6";  Value = 60 }
        @{ Cell = "F17"; Value = 4430 }
        @{ Cell = "F20"; Value = 482 }
        @{ Cell = "F23"; Value = 13 }
        @{ Cell = "F24"; Value = 1108 }
        @{ Cell = "F25"; Value = 2002 }
        @{ Cell = "F30"; Value = 2118 }
        @{ Cell = "F31"; Value = 77 }
        @{ Cell = "F34"; Value = 150 }
        @{ Cell = "F36"; Value = 35 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates[$sheetName]) {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
